# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and
# "全部类型" worksheets to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$ws1.Range("F9").Value = 52
$ws1.Range("F10").Value = 14
$ws1.Range("F12").Value = 2049
$ws1.Range("F15").Value = 1331
$ws1.Range("F16").Value = 466
$ws1.Range("F17").Value = 21
$ws1.Range("F18").Value = 291
$ws1.Range("F19").Value = 209
$ws1.Range("F22").Value = 37
$ws1.Range("F25").Value = 9
$ws1.Range("F26").Value = 1107
$ws1.Range("F27").Value = 7
$ws1.Range("F28").Value = 333
$ws1.Range("F30").Value = 269
$ws1.Range("G30").Value = 55
$ws1.Range("F31").Value = 320

# --- 全部类型 (sheet4) ---
$ws4.Range("F10").Value = 52
$ws4.Range("F11").Value = 14
$ws4.Range("F13").Value = 2049
$ws4.Range("F16").Value = 1331
$ws4.Range("F17").Value = 466
$ws4.Range("F18").Value = 21
$ws4.Range("F19").Value = 291
$ws4.Range("F20").Value = 209
$ws4.Range("F23").Value = 37
$ws4.Range("F26").Value = 9
$ws4.Range("F27").Value = 1107
$ws4.Range("F28").Value = 7
$ws4.Range("F29").Value = 333
$ws4.Range("F31").Value = 269
$ws4.Range("G31").Value = 55
$ws4.Range("F32").Value = 320
